# Auto-generated edit script applying numeric updates to Leve profit tables
# across sheets ALC, ARM, BSM, GSM, LTW, WVR (scheduled runner data refresh).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 377.54
$ws.Range("J17").Value = 377.54
$ws.Range("L17").Value = 1132.62
$ws.Range("N17").Value = -1468.62
$ws.Range("H19").Value = 4358.2856
$ws.Range("I19").Value = 4170.3335
$ws.Range("J19").Value = 4499.25
$ws.Range("K19").Value = 4170.3335
$ws.Range("L19").Value = 4499.25
$ws.Range("M19").Value = -3995.3335
$ws.Range("N19").Value = -4849.25
$ws.Range("H51").Value = 24807.691
$ws.Range("I51").Value = 10642.857
$ws.Range("J51").Value = 41333.332
$ws.Range("K51").Value = 10642.857
$ws.Range("L51").Value = 41333.332
$ws.Range("M51").Value = -10158.857
$ws.Range("N51").Value = -42301.332
$ws.Range("H70").Value = 4133.303
$ws.Range("I70").Value = 3118.182
$ws.Range("J70").Value = 6163.5454
$ws.Range("K70").Value = 9354.545999999998
$ws.Range("L70").Value = 18490.6362
$ws.Range("M70").Value = -9084.545999999998
$ws.Range("N70").Value = -19030.6362
$ws.Range("H73").Value = 4133.303
$ws.Range("I73").Value = 3118.182
$ws.Range("J73").Value = 6163.5454
$ws.Range("K73").Value = 9354.545999999998
$ws.Range("L73").Value = 18490.6362
$ws.Range("M73").Value = -8418.545999999998
$ws.Range("N73").Value = -20362.6362
$ws.Range("H86").Value = 950
$ws.Range("I86").Value = 950
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 950
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = 173
$ws.Range("N86").ClearContents()
$ws.Range("H88").Value = 834.125
$ws.Range("J88").Value = 734.6
$ws.Range("L88").Value = 734.6
$ws.Range("N88").Value = -1546.6
$ws.Range("H89").Value = 950
$ws.Range("I89").Value = 950
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 4750
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = 866
$ws.Range("N89").ClearContents()
$ws.Range("H91").Value = 834.125
$ws.Range("J91").Value = 734.6
$ws.Range("L91").Value = 734.6
$ws.Range("N91").Value = -3542.6
$ws.Range("H107").Value = 26317376
$ws.Range("I107").Value = 33334106
$ws.Range("K107").Value = 33334106
$ws.Range("M107").Value = -33332186
$ws.Range("H113").Value = 156317
$ws.Range("I113").Value = 2452.25
$ws.Range("J113").Value = 224701.33
$ws.Range("K113").Value = 2452.25
$ws.Range("L113").Value = 224701.33
$ws.Range("M113").Value = 801.75
$ws.Range("N113").Value = -231209.33
$ws.Range("H137").Value = 10881735
$ws.Range("I137").Value = 31258274
$ws.Range("K137").Value = 93774822
$ws.Range("M137").Value = -93772272

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3682210
$ws.Range("I2").Value = 5412913.5
$ws.Range("J2").Value = 4465.5
$ws.Range("K2").Value = 5412913.5
$ws.Range("L2").Value = 4465.5
$ws.Range("M2").Value = -5412800.5
$ws.Range("N2").Value = -4691.5
$ws.Range("H32").Value = 21434.1
$ws.Range("I32").Value = 21434.1
$ws.Range("K32").Value = 21434.1
$ws.Range("M32").Value = -21147.1
$ws.Range("H63").Value = 8325.3125
$ws.Range("J63").Value = 11240.477
$ws.Range("L63").Value = 11240.477
$ws.Range("N63").Value = -12612.477
$ws.Range("H66").Value = 8325.3125
$ws.Range("J66").Value = 11240.477
$ws.Range("L66").Value = 56202.385
$ws.Range("N66").Value = -63066.385
$ws.Range("H110").Value = 7813835.5
$ws.Range("I110").Value = 9616105
$ws.Range("K110").Value = 9616105
$ws.Range("M110").Value = -9614060
$ws.Range("H116").Value = 3682210
$ws.Range("I116").Value = 5412913.5
$ws.Range("J116").Value = 4465.5
$ws.Range("K116").Value = 5412913.5
$ws.Range("L116").Value = 4465.5
$ws.Range("M116").Value = -5410619.5
$ws.Range("N116").Value = -9053.5
$ws.Range("H117").Value = 225050000
$ws.Range("J117").Value = 225050000
$ws.Range("L117").Value = 225050000
$ws.Range("N117").Value = -225059178
$ws.Range("H119").Value = 50000
$ws.Range("J119").Value = 50000
$ws.Range("L119").Value = 50000
$ws.Range("N119").Value = -59676
$ws.Range("H122").Value = 3171.375
$ws.Range("I122").Value = 3171.375
$ws.Range("K122").Value = 9514.125
$ws.Range("M122").Value = -7064.125

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3682210
$ws.Range("I3").Value = 5412913.5
$ws.Range("J3").Value = 4465.5
$ws.Range("K3").Value = 5412913.5
$ws.Range("L3").Value = 4465.5
$ws.Range("M3").Value = -5412799.5
$ws.Range("N3").Value = -4693.5
$ws.Range("H99").Value = 4158.8423
$ws.Range("I99").Value = 3971.923
$ws.Range("K99").Value = 3971.923
$ws.Range("M99").Value = -2473.923

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3917.3684
$ws.Range("I80").Value = 3782.2727
$ws.Range("J80").Value = 4103.125
$ws.Range("K80").Value = 3782.2727
$ws.Range("L80").Value = 4103.125
$ws.Range("M80").Value = -2784.2727
$ws.Range("N80").Value = -6099.125
$ws.Range("H83").Value = 3917.3684
$ws.Range("I83").Value = 3782.2727
$ws.Range("J83").Value = 4103.125
$ws.Range("K83").Value = 18911.3635
$ws.Range("L83").Value = 20515.625
$ws.Range("M83").Value = -13919.3635
$ws.Range("N83").Value = -30499.625
$ws.Range("H113").Value = 12721.909
$ws.Range("I113").Value = 13770.941
$ws.Range("K113").Value = 13770.941
$ws.Range("M113").Value = -11600.941
$ws.Range("H122").Value = 9719.743
$ws.Range("I122").Value = 12380.777
$ws.Range("J122").Value = 3732.4167
$ws.Range("K122").Value = 37142.331
$ws.Range("L122").Value = 11197.2501
$ws.Range("M122").Value = -34692.331
$ws.Range("N122").Value = -16097.2501

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2399.6667
$ws.Range("I7").Value = 2399.6667
$ws.Range("K7").Value = 2399.6667
$ws.Range("M7").Value = -2287.6667
$ws.Range("H16").Value = 1343.25
$ws.Range("I16").Value = 1343.25
$ws.Range("K16").Value = 1343.25
$ws.Range("M16").Value = -1173.25
$ws.Range("H21").Value = 6199.6
$ws.Range("J21").Value = 7999
$ws.Range("L21").Value = 7999
$ws.Range("N21").Value = -8347
$ws.Range("H24").Value = 4797
$ws.Range("J24").Value = 4797
$ws.Range("L24").Value = 4797
$ws.Range("N24").Value = -5483
$ws.Range("H25").Value = 11126.75
$ws.Range("I25").Value = 11335.667
$ws.Range("J25").Value = 10500
$ws.Range("K25").Value = 11335.667
$ws.Range("L25").Value = 10500
$ws.Range("M25").Value = -11105.667
$ws.Range("N25").Value = -10960
$ws.Range("H40").Value = 83347170
$ws.Range("I40").Value = 100013800
$ws.Range("K40").Value = 100013800
$ws.Range("M40").Value = -100013664
$ws.Range("H126").Value = 2399.6667
$ws.Range("I126").Value = 2399.6667
$ws.Range("K126").Value = 7199.000100000001
$ws.Range("M126").Value = -4729.000100000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 3000
$ws.Range("I122").Value = 3000
$ws.Range("K122").Value = 9000
$ws.Range("M122").Value = -6550
